# Weekly refresh of the "Hortaliza, Agrícola del Norte S.A. de Arica - Papa" sheet:
# insert two new daily price records (rows 51-52) and push the rest of the
# historical rows down by two (51-102 -> 53-104).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at position 51, shifting existing rows 51:102 down to 53:104.
$ws.Range("A51:A52").EntireRow.Insert()

# New row 51: Red Lady, Región de O'Higgins
$ws.Cells.Item(51, 1).Value = 1
$ws.Cells.Item(51, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(51, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(51, 4).Value = 44894
$ws.Cells.Item(51, 5).Value = 15
$ws.Cells.Item(51, 6).Value = 100114001
$ws.Cells.Item(51, 7).Value = "Papa"
$ws.Cells.Item(51, 8).Value = "Red Lady"
$ws.Cells.Item(51, 9).Value = "1a (cosecha)"
$ws.Cells.Item(51, 10).Value = 1000
$ws.Cells.Item(51, 11).Value = 15000
$ws.Cells.Item(51, 12).Value = 16000
$ws.Cells.Item(51, 13).Value = 15400
$ws.Cells.Item(51, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(51, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(51, 16).Value = 616
$ws.Cells.Item(51, 17).Value = 25
$ws.Cells.Item(51, 18).Value = "Hortaliza"

# New row 52: Rodeo, Región de O'Higgins
$ws.Cells.Item(52, 1).Value = 1
$ws.Cells.Item(52, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(52, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(52, 4).Value = 44894
$ws.Cells.Item(52, 5).Value = 15
$ws.Cells.Item(52, 6).Value = 100114001
$ws.Cells.Item(52, 7).Value = "Papa"
$ws.Cells.Item(52, 8).Value = "Rodeo"
$ws.Cells.Item(52, 9).Value = "1a (cosecha)"
$ws.Cells.Item(52, 10).Value = 800
$ws.Cells.Item(52, 11).Value = 14000
$ws.Cells.Item(52, 12).Value = 15000
$ws.Cells.Item(52, 13).Value = 14500
$ws.Cells.Item(52, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(52, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(52, 16).Value = 580
$ws.Cells.Item(52, 17).Value = 25
$ws.Cells.Item(52, 18).Value = "Hortaliza"
